$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has data in columns A:C (rows 1-20).
# We need to insert a new first column so that:
#   old A (category names)     -> becomes column B
#   old B (PercActivations)    -> becomes column C
#   old C (PercSegmentAreas)   -> becomes column D
# and populate the new column A with a numeric index (0-based row number)
# plus a new header "segments" in B1.

$lastRow = 20

# Insert a new column before column A; this shifts everything right by one
# column and also carries the old column A's direct formatting (style index 1)
# into the new column B.
$ws.Columns("A").Insert()

# --- Header row ---
$ws.Range("B1").Value = "segments"
# Give B1 the same header formatting (bold / centered / bordered) as the
# neighboring header cells (C1 / D1 already have it after the shift).
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# --- Data rows: fill column A with 0-based sequential numbers ---
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# Give column A (the new numeric index column) the bordered/centered style
# that used to belong to the label column, copying it from B2's current
# (inherited) formatting before we strip it.
$ws.Range("B2").Copy()
$ws.Range("A2:A" + $lastRow).PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# The label column (now B) should no longer carry that bordered/bold style -
# only the plain "Normal" formatting, matching the other plain data columns.
$ws.Range("B2:B" + $lastRow).ClearFormats()

$excel.CutCopyMode = 0
